$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (A_AGENT,D_OUTLET,D_ADDRESS,D_CUSTOMER,B_AGENT_OWNER / AddressDetail / AddressDetailID)
# is duplicated as a new row 12 with identical values (bug-fix commit just re-logs the
# same change-row again after a column change).

# A12: same date value/format as A11 (numeric date serial, custom yyyy-mm-dd format).
$ws.Cells.Item(12, 1).Value = 43686
$ws.Cells.Item(12, 1).NumberFormat = "yyyy-mm-dd"

# B12:N12: copy row 11's values verbatim (keeps text-vs-number/boolean typing identical,
# e.g. "10"/"FALSE" stored as text, N as a real boolean) without carrying over any
# explicit cell-level formatting — PasteSpecial(values) writes the values, and resetting
# the style back to "Normal" drops the cell-level style so the columns' own style keeps
# driving the look, exactly like the untouched row 11 cells.
$ws.Range("B11:N11").Copy()
$ws.Range("B12").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("B12:N12").Style = "Normal"
